# Updates cryptos list price (column D) and 1h volume-change (column E) figures
# on the active worksheet, per the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain, unstyled data cell used as a style template: a couple of the new price
# strings (e.g. "1.00", "151.20") look like plain numbers. Assigning them straight
# to .Value would make Excel re-interpret them as numeric and drop formatting such
# as trailing zeros, so those are entered with a leading "'" (forces text) and the
# cell style is then reset to this template to avoid picking up a quote-prefix style.
$normalStyle = $ws.Range("D9").Style

# Row 2: Bitcoin
$ws.Range("D2").Value = "63.783.21"
$ws.Range("E2").Value = "  +1.19%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "2.625.69"
$ws.Range("E3").Value = "  +0.87%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  -0.15%  "

# Row 5: BNB
$ws.Range("D5").Value = "'597.24"
$ws.Range("D5").Style = $normalStyle
$ws.Range("E5").Value = "  -1.03%  "

# Row 6: Solana
$ws.Range("D6").Value = "'151.20"
$ws.Range("D6").Style = $normalStyle
$ws.Range("E6").Value = "  +4.12%  "

# Row 7: USDC
$ws.Range("E7").Value = "  -0.13%  "

# Row 8: XRP
$ws.Range("E8").Value = "  +1.17%  "

# Row 9: Dogecoin
$ws.Range("E9").Value = "  +1.53%  "

# Row 10: Toncoin
$ws.Range("E10").Value = "  +3.22%  "

# Row 11: Cardano
$ws.Range("D11").Value = "'0.386"
$ws.Range("D11").Style = $normalStyle
$ws.Range("E11").Value = "  +4.65%  "

# Row 12: TRON
$ws.Range("E12").Value = "  -0.88%  "

# Row 13: Avalanche
$ws.Range("D13").Value = "'27.93"
$ws.Range("D13").Style = $normalStyle
$ws.Range("E13").Value = "  +2.80%  "

# Row 14: WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "3.096.75"
$ws.Range("E14").Value = "  +0.57%  "

# Row 15: WrappedBTC
$ws.Range("D15").Value = "63.617.39"
$ws.Range("E15").Value = "  +1.09%  "

# Row 16: ShibaInu
$ws.Range("E16").Value = "  +5.03%  "

# Row 17: WrappedEther
$ws.Range("D17").Value = "2.630.57"
$ws.Range("E17").Value = "  +0.59%  "

# Row 18: Chainlink
$ws.Range("D18").Value = "'12.37"
$ws.Range("D18").Style = $normalStyle
$ws.Range("E18").Value = "  +8.26%  "

# Row 19: Polkadot
$ws.Range("D19").Value = "'4.72"
$ws.Range("D19").Style = $normalStyle
$ws.Range("E19").Value = "  +4.48%  "

# Row 20: BitcoinCash
$ws.Range("D20").Value = "'348.13"
$ws.Range("D20").Style = $normalStyle
$ws.Range("E20").Value = "  +2.15%  "

# Row 21: Uniswap
$ws.Range("D21").Value = "'6.89"
$ws.Range("D21").Style = $normalStyle
$ws.Range("E21").Value = "  +1.07%  "

# Row 22: Dai
$ws.Range("D22").Value = "'1.00"
$ws.Range("D22").Style = $normalStyle
$ws.Range("E22").Value = "  +0.04%  "

# Row 23: LEO
$ws.Range("E23").Value = "  +0.83%  "

# Row 24: Litecoin
$ws.Range("D24").Value = "'66.91"
$ws.Range("D24").Style = $normalStyle
$ws.Range("E24").Value = "  +0.81%  "

# Row 25: SuiNetwork
$ws.Range("D25").Value = "'1.72"
$ws.Range("D25").Style = $normalStyle
$ws.Range("E25").Value = "  +9.74%  "

# Row 26: InternetComputer(DFINITY)
$ws.Range("D26").Value = "'9.33"
$ws.Range("D26").Style = $normalStyle
$ws.Range("E26").Value = "  +3.55%  "

# Row 27: Fetch.AI
$ws.Range("E27").Value = "  -0.21%  "

# Row 28: Bittensor
$ws.Range("D28").Value = "'560.15"
$ws.Range("D28").Style = $normalStyle
$ws.Range("E28").Value = "  +1.62%  "

# Row 29: Aptos
$ws.Range("E29").Value = "  +4.72%  "

# Row 30: Kaspa
$ws.Range("E30").Value = "  +0.33%  "

# Row 31: Binance-PegBSC-USD
$ws.Range("E31").Value = "  -0.45%  "

# Row 32: PancakeSwap
$ws.Range("E32").Value = "  +1.04%  "

# Row 33: PEPE
$ws.Range("D33").Value = "0.0$([char]0x2083)0853"
$ws.Range("E33").Value = "  +1.93%  "

# Row 34: ImmutableX
$ws.Range("D34").Value = "'1.76"
$ws.Range("D34").Style = $normalStyle
$ws.Range("E34").Value = "  +1.36%  "

# Row 35: NEARProtocol
$ws.Range("D35").Value = "'5.30"
$ws.Range("D35").Style = $normalStyle
$ws.Range("E35").Value = "  +3.40%  "

# Row 36: Monero
$ws.Range("D36").Value = "'168.11"
$ws.Range("D36").Style = $normalStyle
$ws.Range("E36").Value = "  +0.49%  "

# Row 37: PolygonEcosystemToken
$ws.Range("D37").Value = "'0.416"
$ws.Range("D37").Style = $normalStyle
$ws.Range("E37").Value = "  +3.91%  "

# Row 38: FirstDigitalUSD
$ws.Range("D38").Value = "'1.00"
$ws.Range("D38").Style = $normalStyle
$ws.Range("E38").Value = "  -0.12%  "

# Row 39: EthereumClassic
$ws.Range("D39").Value = "'19.59"
$ws.Range("D39").Style = $normalStyle
$ws.Range("E39").Value = "  +3.41%  "

# Row 40: Stacks
$ws.Range("E40").Value = "  +1.43%  "

# Row 41: USDe
$ws.Range("E41").Value = "  -0.05%  "

# Row 42: Aave
$ws.Range("D42").Value = "'167.33"
$ws.Range("D42").Style = $normalStyle
$ws.Range("E42").Value = "  +1.63%  "

# Row 43: OKB
$ws.Range("D43").Value = "'39.69"
$ws.Range("D43").Style = $normalStyle
$ws.Range("E43").Value = "  +0.29%  "

# Row 44: Filecoin
$ws.Range("D44").Value = "'3.95"
$ws.Range("D44").Style = $normalStyle
$ws.Range("E44").Value = "  +5.65%  "

# Row 45: Hedera
$ws.Range("D45").Value = "'0.0592"
$ws.Range("D45").Style = $normalStyle
$ws.Range("E45").Value = "  +5.51%  "

# Row 46: InjectiveProtocol
$ws.Range("D46").Value = "'21.95"
$ws.Range("D46").Style = $normalStyle
$ws.Range("E46").Value = "  +1.26%  "

# Row 47: Mantle
$ws.Range("D47").Value = "'0.634"
$ws.Range("D47").Style = $normalStyle
$ws.Range("E47").Value = "  +2.09%  "

# Row 48: VeChain
$ws.Range("E48").Value = "  +3.45%  "

# Row 49: dogwifhat
$ws.Range("D49").Value = "'2.02"
$ws.Range("D49").Style = $normalStyle
$ws.Range("E49").Value = "  +6.71%  "

# Row 50: BabyDogeCoin
$ws.Range("E50").Value = "  +27.24%  "

# Row 51: Stellar
$ws.Range("D51").Value = "'0.0968"
$ws.Range("D51").Style = $normalStyle
$ws.Range("E51").Value = "  +1.57%  "
